$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) is stored as plain text in the source data (values
# like "100.60" or "42.723.07" are not real numbers). Force NumberFormat to
# text ("@") before writing so Excel does not auto-convert numeric-looking
# strings to actual numbers (which would also silently drop significant
# trailing zeros, e.g. "100.60" -> 100.6).
$priceCells = @("D2", "D3", "D5", "D6", "D10", "D12", "D14", "D15", "D16", "D17", "D18", "D20", "D21", "D23", "D24", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D40", "D41", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "42.723.07"
$ws.Range("E2").Value = "  -1.33%  "

$ws.Range("D3").Value = "2.238.76"
$ws.Range("E3").Value = "  -1.61%  "

$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").Value = "114.27"
$ws.Range("E5").Value = "  +2.31%  "

$ws.Range("D6").Value = "278.04"
$ws.Range("E6").Value = "  +5.43%  "

$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("E9").Value = "  +0.71%  "

$ws.Range("D10").Value = "46.35"
$ws.Range("E10").Value = "  -0.72%  "

$ws.Range("E11").Value = "  -0.36%  "

$ws.Range("D12").Value = "9.05"
$ws.Range("E12").Value = "  -1.46%  "

$ws.Range("E13").Value = "  -2.92%  "

$ws.Range("D14").Value = "15.31"
$ws.Range("E14").Value = "  -0.33%  "

$ws.Range("D15").Value = "0.874"
$ws.Range("E15").Value = "  +1.64%  "

$ws.Range("D16").Value = "2.577.08"
$ws.Range("E16").Value = "  -1.63%  "

$ws.Range("D17").Value = "2.239.98"
$ws.Range("E17").Value = "  -1.89%  "

$ws.Range("D18").Value = "42.990.20"
$ws.Range("E18").Value = "  -0.45%  "

$ws.Range("E19").Value = "  -0.62%  "

$ws.Range("D20").Value = "6.78"
$ws.Range("E20").Value = "  +0.77%  "

$ws.Range("D21").Value = "72.19"
$ws.Range("E21").Value = "  +0.62%  "

$ws.Range("E22").Value = "  -4.31%  "

$ws.Range("D23").Value = "3.03"
$ws.Range("E23").Value = "  +6.90%  "

$ws.Range("D24").Value = "231.99"
$ws.Range("E24").Value = "  -0.86%  "

$ws.Range("E25").Value = "  -0.88%  "

$ws.Range("D26").Value = "12.16"
$ws.Range("E26").Value = "  +7.61%  "

$ws.Range("E27").Value = "  -1.05%  "

$ws.Range("D28").Value = "40.41"
$ws.Range("E28").Value = "  +0.16%  "

$ws.Range("D29").Value = "2.25"
$ws.Range("E29").Value = "  +0.08%  "

$ws.Range("D30").Value = "3.26"
$ws.Range("E30").Value = "  -2.33%  "

$ws.Range("D31").Value = "173.21"
$ws.Range("E31").Value = "  +0.26%  "

$ws.Range("D32").Value = "21.08"
$ws.Range("E32").Value = "  -1.54%  "

$ws.Range("D33").Value = "0.0895"
$ws.Range("E33").Value = "  -0.50%  "

$ws.Range("D34").Value = "5.57"
$ws.Range("E34").Value = "  -0.84%  "

$ws.Range("D35").Value = "4.39"
$ws.Range("E35").Value = "  +9.31%  "

$ws.Range("E36").Value = "  -0.20%  "

$ws.Range("D37").Value = "0.0374"
$ws.Range("E37").Value = "  +2.00%  "

$ws.Range("D38").Value = "4.64"
$ws.Range("E38").Value = "  +0.54%  "

$ws.Range("E39").Value = "  +2.55%  "

$ws.Range("D40").Value = "2.56"
$ws.Range("E40").Value = "  -1.45%  "

$ws.Range("D41").Value = "71.21"
$ws.Range("E41").Value = "  -6.72%  "

$ws.Range("E42").Value = "  -1.76%  "

$ws.Range("D43").Value = "13.19"
$ws.Range("E43").Value = "  -5.18%  "

$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "1.34"
$ws.Range("E45").Value = "  -2.53%  "

$ws.Range("B46").Value = "THORChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D46").Value = "5.66"
$ws.Range("E46").Value = "  -6.54%  "

$ws.Range("D47").Value = "1.28"
$ws.Range("E47").Value = "  +2.05%  "

$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "8.44"
$ws.Range("E48").Value = "  -1.34%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0991"
$ws.Range("E49").Value = "  -0.45%  "

$ws.Range("D50").Value = "100.60"
$ws.Range("E50").Value = "  -2.19%  "

$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").Value = "0.642"
$ws.Range("E51").Value = "  +7.98%  "
